$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.280371333333333
$ws.Range("H2").Value = 24.841114
$ws.Range("I2").Value = 0.2946400644635011
$ws.Range("J2").Value = 0.3116548779253407
$ws.Range("M2").Value = 0.6598136666666666
$ws.Range("N2").Value = 1.979441
$ws.Range("O2").Value = 0.007704735356083927
$ws.Range("P2").Value = 0.008484678519943686
$ws.Range("Q2").Value = 5.463502170808222
$ws.Range("R2").Value = 49.171519537274
$ws.Range("S2").Value = 0.002270123721990784
$ws.Range("T2").Value = 0.00264429144836881
$ws.Range("G3").Value = 8.280371333333333
$ws.Range("H3").Value = 24.841114
$ws.Range("I3").Value = 0.2946400644635011
$ws.Range("J3").Value = 0.3116548779253407
$ws.Range("O3").Value = 0.7130079175842846
$ws.Range("P3").Value = 0.7851850431306702
$ws.Range("Q3").Value = 505.6007929524991
$ws.Range("R3").Value = 4550.407136572492
$ws.Range("S3").Value = 0.2100806988000203
$ws.Range("T3").Value = 0.2447067487656924
$ws.Range("G4").Value = 8.280371333333333
$ws.Range("H4").Value = 24.841114
$ws.Range("I4").Value = 0.2946400644635011
$ws.Range("J4").Value = 0.3116548779253407
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1563486666666667
$ws.Range("N4").Value = 0.469046
$ws.Range("O4").Value = 0.001825704984300993
$ws.Range("P4").Value = 0.002010519394650058
$ws.Range("Q4").Value = 1.294625017471555
$ws.Range("R4").Value = 11.651625157244
$ws.Range("S4").Value = 0.0005379258342657798
$ws.Range("T4").Value = 0.000626588176506194
$ws.Range("G5").Value = 8.280371333333333
$ws.Range("H5").Value = 24.841114
$ws.Range("I5").Value = 0.2946400644635011
$ws.Range("J5").Value = 0.3116548779253407
$ws.Range("M5").Value = 23.6163295
$ws.Range("N5").Value = 47.232659
$ws.Range("O5").Value = 0.2757711427815902
$ws.Range("P5").Value = 0.2024581319964196
$ws.Range("Q5").Value = 195.5519777903543
$ws.Range("R5").Value = 1173.311866742126
$ws.Range("S5").Value = 0.08125322728634109
$ws.Range("T5").Value = 0.06309706441233667
$ws.Range("G6").Value = 8.280371333333333
$ws.Range("H6").Value = 24.841114
$ws.Range("I6").Value = 0.2946400644635011
$ws.Range("J6").Value = 0.3116548779253407
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.14477
$ws.Range("N6").Value = 0.43431
$ws.Range("O6").Value = 0.00169049929374041
$ws.Range("P6").Value = 0.001861626958316384
$ws.Range("Q6").Value = 1.198749357926667
$ws.Range("R6").Value = 10.78874422134
$ws.Range("S6").Value = 0.0004980888208831775
$ws.Range("T6").Value = 0.0005801851224366162
$ws.Range("I7").Value = 0.405746032520008
$ws.Range("J7").Value = 0.4291769704298953
$ws.Range("M7").Value = 0.6598136666666666
$ws.Range("N7").Value = 1.979441
$ws.Range("O7").Value = 0.007704735356083927
$ws.Range("P7").Value = 0.008484678519943686
$ws.Range("Q7").Value = 7.523736914415776
$ws.Range("R7").Value = 67.71363222974199
$ws.Range("S7").Value = 0.003126165802347685
$ws.Range("T7").Value = 0.003641428622261039
$ws.Range("I8").Value = 0.405746032520008
$ws.Range("J8").Value = 0.4291769704298953
$ws.Range("O8").Value = 0.7130079175842846
$ws.Range("P8").Value = 0.7851850431306702
$ws.Range("S8").Value = 0.2893001337151763
$ws.Range("T8").Value = 0.3369833380376877
$ws.Range("I9").Value = 0.405746032520008
$ws.Range("J9").Value = 0.4291769704298953
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1563486666666667
$ws.Range("N9").Value = 0.469046
$ws.Range("O9").Value = 0.001825704984300993
$ws.Range("P9").Value = 0.002010519394650058
$ws.Range("Q9").Value = 1.782815807472444
$ws.Range("R9").Value = 16.045342267252
$ws.Range("S9").Value = 0.0007407725539321314
$ws.Range("T9").Value = 0.000862868622786459
$ws.Range("I10").Value = 0.405746032520008
$ws.Range("J10").Value = 0.4291769704298953
$ws.Range("M10").Value = 23.6163295
$ws.Range("N10").Value = 47.232659
$ws.Range("O10").Value = 0.2757711427815902
$ws.Range("P10").Value = 0.2024581319964196
$ws.Range("Q10").Value = 269.2927700934096
$ws.Range("R10").Value = 1615.756620560458
$ws.Range("S10").Value = 0.1118930470671388
$ws.Range("T10").Value = 0.08689036772911921
$ws.Range("I11").Value = 0.405746032520008
$ws.Range("J11").Value = 0.4291769704298953
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.14477
$ws.Range("N11").Value = 0.43431
$ws.Range("O11").Value = 0.00169049929374041
$ws.Range("P11").Value = 0.001861626958316384
$ws.Range("Q11").Value = 1.650786347913333
$ws.Range("R11").Value = 14.85707713122
$ws.Range("S11").Value = 0.0006859133814130469
$ws.Range("T11").Value = 0.0007989674180408469
$ws.Range("G12").Value = 1.864050333333333
$ws.Range("H12").Value = 5.592150999999999
$ws.Range("I12").Value = 0.06632841551025578
$ws.Range("J12").Value = 0.07015873512134246
$ws.Range("M12").Value = 0.6598136666666666
$ws.Range("N12").Value = 1.979441
$ws.Range("O12").Value = 0.007704735356083927
$ws.Range("P12").Value = 0.008484678519943686
$ws.Range("Q12").Value = 1.229925885287889
$ws.Range("R12").Value = 11.069332967591
$ws.Range("S12").Value = 0.0005110428880948933
$ws.Range("T12").Value = 0.000595274312870473
$ws.Range("G13").Value = 1.864050333333333
$ws.Range("H13").Value = 5.592150999999999
$ws.Range("I13").Value = 0.06632841551025578
$ws.Range("J13").Value = 0.07015873512134246
$ws.Range("O13").Value = 0.7130079175842846
$ws.Range("P13").Value = 0.7851850431306702
$ws.Range("Q13").Value = 113.8192103586864
$ws.Range("R13").Value = 1024.372893228178
$ws.Range("S13").Value = 0.04729268541963264
$ws.Range("T13").Value = 0.05508758946224455
$ws.Range("G14").Value = 1.864050333333333
$ws.Range("H14").Value = 5.592150999999999
$ws.Range("I14").Value = 0.06632841551025578
$ws.Range("J14").Value = 0.07015873512134246
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1563486666666667
$ws.Range("N14").Value = 0.469046
$ws.Range("O14").Value = 0.001825704984300993
$ws.Range("P14").Value = 0.002010519394650058
$ws.Range("Q14").Value = 0.2914417842162222
$ws.Range("R14").Value = 2.622976057946
$ws.Range("S14").Value = 0.0001210961187978613
$ws.Range("T14").Value = 0.0001410554976655752
$ws.Range("G15").Value = 1.864050333333333
$ws.Range("H15").Value = 5.592150999999999
$ws.Range("I15").Value = 0.06632841551025578
$ws.Range("J15").Value = 0.07015873512134246
$ws.Range("M15").Value = 23.6163295
$ws.Range("N15").Value = 47.232659
$ws.Range("O15").Value = 0.2757711427815902
$ws.Range("P15").Value = 0.2024581319964196
$ws.Range("Q15").Value = 44.02202687658482
$ws.Range("R15").Value = 264.132161259509
$ws.Range("S15").Value = 0.01829146294415539
$ws.Range("T15").Value = 0.01420420645589859
$ws.Range("G16").Value = 1.864050333333333
$ws.Range("H16").Value = 5.592150999999999
$ws.Range("I16").Value = 0.06632841551025578
$ws.Range("J16").Value = 0.07015873512134246
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.14477
$ws.Range("N16").Value = 0.43431
$ws.Range("O16").Value = 0.00169049929374041
$ws.Range("P16").Value = 0.001861626958316384
$ws.Range("Q16").Value = 0.2698585667566666
$ws.Range("R16").Value = 2.42872710081
$ws.Range("S16").Value = 0.0001121281395750079
$ws.Range("T16").Value = 0.0001306093926632697
$ws.Range("G17").Value = 4.6029105
$ws.Range("H17").Value = 9.205821
$ws.Range("I17").Value = 0.1637851482553954
$ws.Range("J17").Value = 0.1154955860658076
$ws.Range("M17").Value = 0.6598136666666666
$ws.Range("N17").Value = 1.979441
$ws.Range("O17").Value = 0.007704735356083927
$ws.Range("P17").Value = 0.008484678519943686
$ws.Range("Q17").Value = 3.0370632543435
$ws.Range("R17").Value = 18.222379526061
$ws.Range("S17").Value = 0.001261921222564793
$ws.Range("T17").Value = 0.000979942918240865
$ws.Range("G18").Value = 4.6029105
$ws.Range("H18").Value = 9.205821
$ws.Range("I18").Value = 0.1637851482553954
$ws.Range("J18").Value = 0.1154955860658076
$ws.Range("O18").Value = 0.7130079175842846
$ws.Range("P18").Value = 0.7851850431306702
$ws.Range("Q18").Value = 281.054448527073
$ws.Range("R18").Value = 1686.326691162438
$ws.Range("S18").Value = 0.1167801074888128
$ws.Range("T18").Value = 0.09068540672648318
$ws.Range("G19").Value = 4.6029105
$ws.Range("H19").Value = 9.205821
$ws.Range("I19").Value = 0.1637851482553954
$ws.Range("J19").Value = 0.1154955860658076
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.1563486666666667
$ws.Range("N19").Value = 0.469046
$ws.Range("O19").Value = 0.001825704984300993
$ws.Range("P19").Value = 0.002010519394650058
$ws.Range("Q19").Value = 0.719658919461
$ws.Range("R19").Value = 4.317953516766
$ws.Range("S19").Value = 0.0002990233615243525
$ws.Range("T19").Value = 0.0002322061157817812
$ws.Range("G20").Value = 4.6029105
$ws.Range("H20").Value = 9.205821
$ws.Range("I20").Value = 0.1637851482553954
$ws.Range("J20").Value = 0.1154955860658076
$ws.Range("M20").Value = 23.6163295
$ws.Range("N20").Value = 47.232659
$ws.Range("O20").Value = 0.2757711427815902
$ws.Range("P20").Value = 0.2024581319964196
$ws.Range("Q20").Value = 108.7038510270097
$ws.Range("R20").Value = 434.815404108039
$ws.Range("S20").Value = 0.04516721750504256
$ws.Range("T20").Value = 0.02338302060871512
$ws.Range("G21").Value = 4.6029105
$ws.Range("H21").Value = 9.205821
$ws.Range("I21").Value = 0.1637851482553954
$ws.Range("J21").Value = 0.1154955860658076
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.14477
$ws.Range("N21").Value = 0.43431
$ws.Range("O21").Value = 0.00169049929374041
$ws.Range("P21").Value = 0.001861626958316384
$ws.Range("Q21").Value = 0.6663633530850001
$ws.Range("R21").Value = 3.998180118510001
$ws.Range("S21").Value = 0.0002768786774509143
$ws.Range("T21").Value = 0.0002150096965866576
$ws.Range("G22").Value = 1.953192
$ws.Range("H22").Value = 5.859576
$ws.Range("I22").Value = 0.0695003392508397
$ws.Range("J22").Value = 0.0735138304576138
$ws.Range("M22").Value = 0.6598136666666666
$ws.Range("N22").Value = 1.979441
$ws.Range("O22").Value = 0.007704735356083927
$ws.Range("P22").Value = 0.008484678519943686
$ws.Range("Q22").Value = 1.288742775224
$ws.Range("R22").Value = 11.598684977016
$ws.Range("S22").Value = 0.0005354817210857721
$ws.Range("T22").Value = 0.0006237412182024977
$ws.Range("G23").Value = 1.953192
$ws.Range("H23").Value = 5.859576
$ws.Range("I23").Value = 0.0695003392508397
$ws.Range("J23").Value = 0.0735138304576138
$ws.Range("O23").Value = 0.7130079175842846
$ws.Range("P23").Value = 0.7851850431306702
$ws.Range("Q23").Value = 119.262214728592
$ws.Range("R23").Value = 1073.359932557328
$ws.Range("S23").Value = 0.04955429216064253
$ws.Range("T23").Value = 0.05772196013856226
$ws.Range("G24").Value = 1.953192
$ws.Range("H24").Value = 5.859576
$ws.Range("I24").Value = 0.0695003392508397
$ws.Range("J24").Value = 0.0735138304576138
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.1563486666666667
$ws.Range("N24").Value = 0.469046
$ws.Range("O24").Value = 0.001825704984300993
$ws.Range("P24").Value = 0.002010519394650058
$ws.Range("Q24").Value = 0.305378964944
$ws.Range("R24").Value = 2.748410684496
$ws.Range("S24").Value = 0.000126887115780868
$ws.Range("T24").Value = 0.0001478009819100487
$ws.Range("G25").Value = 1.953192
$ws.Range("H25").Value = 5.859576
$ws.Range("I25").Value = 0.0695003392508397
$ws.Range("J25").Value = 0.0735138304576138
$ws.Range("M25").Value = 23.6163295
$ws.Range("N25").Value = 47.232659
$ws.Range("O25").Value = 0.2757711427815902
$ws.Range("P25").Value = 0.2024581319964196
$ws.Range("Q25").Value = 46.12722584876399
$ws.Range("R25").Value = 276.763355092584
$ws.Range("S25").Value = 0.01916618797891227
$ws.Range("T25").Value = 0.01488347279034999
$ws.Range("G26").Value = 1.953192
$ws.Range("H26").Value = 5.859576
$ws.Range("I26").Value = 0.0695003392508397
$ws.Range("J26").Value = 0.0735138304576138
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.14477
$ws.Range("N26").Value = 0.43431
$ws.Range("O26").Value = 0.00169049929374041
$ws.Range("P26").Value = 0.001861626958316384
$ws.Range("Q26").Value = 0.28276360584
$ws.Range("R26").Value = 2.54487245256
$ws.Range("S26").Value = 0.0001121281395750079
$ws.Range("T26").Value = 0.0001368553285889939
